$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 18699.8
$ws.Range("J21").Value = 18699.8
$ws.Range("L21").Value = 18699.8
$ws.Range("N21").Value = -19635.8
$ws.Range("H23").Value = 18699.8
$ws.Range("J23").Value = 18699.8
$ws.Range("L23").Value = 18699.8
$ws.Range("N23").Value = -19167.8
$ws.Range("H38").Value = 408.0625
$ws.Range("I38").Value = 244.08333
$ws.Range("K38").Value = 732.24999
$ws.Range("M38").Value = -360.24999
$ws.Range("H45").Value = 1758.5
$ws.Range("J45").Value = 2500
$ws.Range("L45").Value = 7500
$ws.Range("N45").Value = -7884
$ws.Range("H58").Value = 1079.7
$ws.Range("I58").Value = 812.93335
$ws.Range("J58").Value = 1880
$ws.Range("K58").Value = 2438.80005
$ws.Range("L58").Value = 5640
$ws.Range("M58").Value = -2288.80005
$ws.Range("N58").Value = -5940
$ws.Range("H87").Value = 25190.592
$ws.Range("J87").Value = 25190.592
$ws.Range("L87").Value = 25190.592
$ws.Range("N87").Value = -27686.592
$ws.Range("H90").Value = 25190.592
$ws.Range("J90").Value = 25190.592
$ws.Range("L90").Value = 75571.776
$ws.Range("N90").Value = -88051.776
$ws.Range("H141").Value = 5405.3447
$ws.Range("I141").Value = 3009.1052
$ws.Range("J141").Value = 9958.200000000001
$ws.Range("K141").Value = 9027.3156
$ws.Range("L141").Value = 29874.6
$ws.Range("M141").Value = -3847.3156
$ws.Range("N141").Value = -40234.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 12524.048
$ws.Range("J23").Value = 9649.950000000001
$ws.Range("L23").Value = 9649.950000000001
$ws.Range("N23").Value = -10167.95
$ws.Range("H44").Value = 29332
$ws.Range("J44").Value = 29332
$ws.Range("L44").Value = 29332
$ws.Range("N44").Value = -30308
$ws.Range("H63").Value = 5109.8887
$ws.Range("I63").Value = 3333.3333
$ws.Range("J63").Value = 5998.1665
$ws.Range("K63").Value = 3333.3333
$ws.Range("L63").Value = 5998.1665
$ws.Range("M63").Value = -2647.3333
$ws.Range("N63").Value = -7370.1665
$ws.Range("H66").Value = 5109.8887
$ws.Range("I66").Value = 3333.3333
$ws.Range("J66").Value = 5998.1665
$ws.Range("K66").Value = 16666.6665
$ws.Range("L66").Value = 29990.8325
$ws.Range("M66").Value = -13234.6665
$ws.Range("N66").Value = -36854.8325
$ws.Range("H74").Value = 1906.6875
$ws.Range("I74").Value = 1288.7778
$ws.Range("J74").Value = 2701.1428
$ws.Range("K74").Value = 1288.7778
$ws.Range("L74").Value = 2701.1428
$ws.Range("M74").Value = -414.7778000000001
$ws.Range("N74").Value = -4449.1428
$ws.Range("H77").Value = 1906.6875
$ws.Range("I77").Value = 1288.7778
$ws.Range("J77").Value = 2701.1428
$ws.Range("K77").Value = 6443.889
$ws.Range("L77").Value = 13505.714
$ws.Range("M77").Value = -2075.889
$ws.Range("N77").Value = -22241.714
$ws.Range("H80").Value = 21665.666
$ws.Range("J80").Value = 21665.666
$ws.Range("L80").Value = 21665.666
$ws.Range("N80").Value = -23661.666
$ws.Range("H83").Value = 21665.666
$ws.Range("J83").Value = 21665.666
$ws.Range("L83").Value = 64996.99800000001
$ws.Range("N83").Value = -74980.99800000001
$ws.Range("H122").Value = 1430.3016
$ws.Range("I122").Value = 1283.9706
$ws.Range("K122").Value = 3851.9118
$ws.Range("M122").Value = -1401.9118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 14897.4
$ws.Range("J35").Value = 17371.75
$ws.Range("L35").Value = 17371.75
$ws.Range("N35").Value = -17991.75
$ws.Range("H82").Value = 11808.385
$ws.Range("J82").Value = 21999
$ws.Range("L82").Value = 21999
$ws.Range("N82").Value = -22765
$ws.Range("H85").Value = 11808.385
$ws.Range("J85").Value = 21999
$ws.Range("L85").Value = 21999
$ws.Range("N85").Value = -24651
$ws.Range("H99").Value = 1468
$ws.Range("I99").Value = 1458.8
$ws.Range("J99").Value = 1483.3334
$ws.Range("K99").Value = 1458.8
$ws.Range("L99").Value = 1483.3334
$ws.Range("M99").Value = 39.20000000000005
$ws.Range("N99").Value = -4479.3334
$ws.Range("H105").Value = 25002854
$ws.Range("I105").Value = 41669656
$ws.Range("K105").Value = 41669656
$ws.Range("M105").Value = -41667909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 60326.668
$ws.Range("J28").Value = 60326.668
$ws.Range("L28").Value = 60326.668
$ws.Range("N28").Value = -60816.668
$ws.Range("H43").Value = 134000
$ws.Range("J43").Value = 134000
$ws.Range("L43").Value = 134000
$ws.Range("N43").Value = -134368
$ws.Range("H50").Value = 9999
$ws.Range("J50").Value = 9999
$ws.Range("L50").Value = 9999
$ws.Range("N50").Value = -11249
$ws.Range("H59").Value = 18499.166
$ws.Range("J59").Value = 18499.166
$ws.Range("L59").Value = 18499.166
$ws.Range("N59").Value = -20789.166
$ws.Range("H60").Value = 10667.667
$ws.Range("J60").Value = 10667.667
$ws.Range("L60").Value = 10667.667
$ws.Range("N60").Value = -11689.667
$ws.Range("H68").Value = 23950.846
$ws.Range("J68").Value = 23950.846
$ws.Range("L68").Value = 23950.846
$ws.Range("N68").Value = -25448.846
$ws.Range("H71").Value = 23950.846
$ws.Range("J71").Value = 23950.846
$ws.Range("L71").Value = 71852.538
$ws.Range("N71").Value = -79340.538
$ws.Range("H74").Value = 25251
$ws.Range("J74").Value = 25251
$ws.Range("L74").Value = 25251
$ws.Range("N74").Value = -26999
$ws.Range("H77").Value = 25251
$ws.Range("J77").Value = 25251
$ws.Range("L77").Value = 75753
$ws.Range("N77").Value = -84489
$ws.Range("H95").Value = 16399.4
$ws.Range("J95").Value = 16399.4
$ws.Range("L95").Value = 16399.4
$ws.Range("N95").Value = -21891.4
$ws.Range("H101").Value = 134000
$ws.Range("J101").Value = 134000
$ws.Range("L101").Value = 134000
$ws.Range("N101").Value = -140490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 447.92682
$ws.Range("I34").Value = 98
$ws.Range("J34").Value = 507.91428
$ws.Range("K34").Value = 294
$ws.Range("L34").Value = 1523.74284
$ws.Range("M34").Value = -210
$ws.Range("N34").Value = -1691.74284
$ws.Range("H39").Value = 1935
$ws.Range("J39").Value = 1935
$ws.Range("L39").Value = 5805
$ws.Range("N39").Value = -6393
$ws.Range("H55").Value = 1738.0952
$ws.Range("J55").Value = 1735
$ws.Range("L55").Value = 5205
$ws.Range("N55").Value = -5559
$ws.Range("H59").Value = 2999.9443
$ws.Range("I59").Value = 2999
$ws.Range("K59").Value = 8997
$ws.Range("M59").Value = -8457
$ws.Range("H111").Value = 6876.6665
$ws.Range("I111").Value = 753.3333
$ws.Range("K111").Value = 2259.9999
$ws.Range("M111").Value = 807.0001000000002
$ws.Range("H126").Value = 3038.7273
$ws.Range("I126").Value = 1715
$ws.Range("J126").Value = 3332.889
$ws.Range("K126").Value = 5145
$ws.Range("L126").Value = 9998.667000000001
$ws.Range("M126").Value = -205
$ws.Range("N126").Value = -19878.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 10866.333
$ws.Range("J43").Value = 15974.5
$ws.Range("L43").Value = 15974.5
$ws.Range("N43").Value = -16276.5
$ws.Range("H46").Value = 5477.619
$ws.Range("J46").Value = 4421
$ws.Range("L46").Value = 4421
$ws.Range("N46").Value = -4733
$ws.Range("H57").Value = 24460
$ws.Range("J57").Value = 24460
$ws.Range("L57").Value = 24460
$ws.Range("N57").Value = -26100
$ws.Range("H80").Value = 31840280
$ws.Range("I80").Value = 46275724
$ws.Range("J80").Value = 82299.8
$ws.Range("K80").Value = 46275724
$ws.Range("L80").Value = 82299.8
$ws.Range("M80").Value = -46274726
$ws.Range("N80").Value = -84295.8
$ws.Range("H83").Value = 31840280
$ws.Range("I83").Value = 46275724
$ws.Range("J83").Value = 82299.8
$ws.Range("K83").Value = 231378620
$ws.Range("L83").Value = 411499
$ws.Range("M83").Value = -231373628
$ws.Range("N83").Value = -421483
$ws.Range("H109").Value = 10190
$ws.Range("J109").Value = 10190
$ws.Range("L109").Value = 10190
$ws.Range("N109").Value = -12270
$ws.Range("H113").Value = 1044.2667
$ws.Range("I113").Value = 827.36365
$ws.Range("J113").Value = 1640.75
$ws.Range("K113").Value = 827.36365
$ws.Range("L113").Value = 1640.75
$ws.Range("M113").Value = 1342.63635
$ws.Range("N113").Value = -5980.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3048.173
$ws.Range("I122").Value = 2914.6487
$ws.Range("J122").Value = 3377.5334
$ws.Range("K122").Value = 8743.946100000001
$ws.Range("L122").Value = 10132.6002
$ws.Range("M122").Value = -6293.946100000001
$ws.Range("N122").Value = -15032.6002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 37833.168
$ws.Range("J56").Value = 37399.8
$ws.Range("L56").Value = 37399.8
$ws.Range("N56").Value = -38827.8
$ws.Range("H81").Value = 4897.8
$ws.Range("I81").Value = 3747
$ws.Range("J81").Value = 7199.4
$ws.Range("K81").Value = 7494
$ws.Range("L81").Value = 14398.8
$ws.Range("M81").Value = -6433
$ws.Range("N81").Value = -16520.8
$ws.Range("H84").Value = 4897.8
$ws.Range("I84").Value = 3747
$ws.Range("J84").Value = 7199.4
$ws.Range("K84").Value = 37470
$ws.Range("L84").Value = 71994
$ws.Range("M84").Value = -32166
$ws.Range("N84").Value = -82602
